$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New arrival rows appended to the "Main Data" sheet (Friday, Jan 13 continues)
# Row 21 - FR2473 from London (STN), Ryanair B38M
$ws.Range("A21").Value = 20
$ws.Range("B21").Value = "Friday, Jan 13"
$ws.Range("C21").Value = "3:30 PM"
$ws.Range("D21").Value = "FR2473"
$ws.Range("E21").Value = "London"
$ws.Range("F21").Value = "(STN)"
$ws.Range("G21").Value = "Ryanair "
$ws.Range("H21").Value = "B38M"
$ws.Range("I21").Value = "(EI-HMW)"
$ws.Range("J21").Value = "3:16 PM"
$ws.Range("K21").Borders.LineStyle = -4142
$ws.Range("L21").Value = "0 hours, -14 minutes"
$ws.Range("M21").Borders.LineStyle = -4142

# Row 22 - W95175 from London (LTN), Wizz Air A321
$ws.Range("A22").Value = 21
$ws.Range("B22").Value = "Friday, Jan 13"
$ws.Range("C22").Value = "6:55 PM"
$ws.Range("D22").Value = "W95175"
$ws.Range("E22").Value = "London"
$ws.Range("F22").Value = "(LTN)"
$ws.Range("G22").Value = "Wizz Air "
$ws.Range("H22").Value = "A321"
$ws.Range("I22").Value = "(G-WUKG)"
$ws.Range("J22").Value = "6:46 PM"
$ws.Range("K22").Borders.LineStyle = -4142
$ws.Range("L22").Value = "0 hours, -9 minutes"
$ws.Range("M22").Borders.LineStyle = -4142

# Row 23 - FR5106 from Dublin (DUB), Ryanair B738
$ws.Range("A23").Value = 22
$ws.Range("B23").Value = "Friday, Jan 13"
$ws.Range("C23").Value = "7:35 PM"
$ws.Range("D23").Value = "FR5106"
$ws.Range("E23").Value = "Dublin"
$ws.Range("F23").Value = "(DUB)"
$ws.Range("G23").Value = "Ryanair "
$ws.Range("H23").Value = "B738"
$ws.Range("I23").Value = "(EI-EMF)"
$ws.Range("J23").Value = "7:11 PM"
$ws.Range("K23").Borders.LineStyle = -4142
$ws.Range("L23").Value = "0 hours, -24 minutes"
$ws.Range("M23").Borders.LineStyle = -4142

# Row 24 - FR9503 from Bristol (BRS), Ryanair B738
$ws.Range("A24").Value = 23
$ws.Range("B24").Value = "Friday, Jan 13"
$ws.Range("C24").Value = "9:25 PM"
$ws.Range("D24").Value = "FR9503"
$ws.Range("E24").Value = "Bristol"
$ws.Range("F24").Value = "(BRS)"
$ws.Range("G24").Value = "Ryanair "
$ws.Range("H24").Value = "B738"
$ws.Range("I24").Value = "(EI-DWH)"
$ws.Range("J24").Value = "9:03 PM"
$ws.Range("K24").Borders.LineStyle = -4142
$ws.Range("L24").Value = "0 hours, -22 minutes"
$ws.Range("M24").Borders.LineStyle = -4142
